$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) and Volume column (E) retain text formatting
# so Excel does not auto-convert numeric-looking strings to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.442.05"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.838.36"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.91%  "
$ws.Range("D5").Value = "243.22"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "0.6263"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.07406"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "0.2933"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "23.33"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "0.07644"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "1.833.50"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "5.010"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "0.6749"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "83.29"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "0.000009326"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "5.892"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "29.410.46"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "2.081.12"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "237.78"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "7.342"
$ws.Range("E23").Value = "  +2.50%  "
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "158.74"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "0.1412"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "8.494"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "17.74"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").Value = "0.06052"
$ws.Range("E29").Value = "  +8.11%  "
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "1.230"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "4.090"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "4.106"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "0.7253"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").Value = "2.612"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").Value = "2.881"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "1.216.24"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "6.294"
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").Value = "0.9114"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "1.996.95"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "101.92"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "65.46"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D49").Value = "9.243"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").Value = "0.4054"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  +2.99%  "

# Rows 47 and 48 swap contents (Mantle <-> BabyDogeCoin) with updated values
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5064"
$ws.Range("E48").Value = "  -1.20%  "
